$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing so that
# numeric-looking strings (e.g. "82.70", "0.9990") keep their exact
# textual representation instead of being coerced to floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.807.85"
$ws.Range("E2").Value = "  +7.80%  "
$ws.Range("D3").Value = "1.743.87"
$ws.Range("E3").Value = "  +4.65%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "335.52"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "0.3743"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  +3.75%  "
$ws.Range("D9").Value = "0.3391"
$ws.Range("E9").Value = "  +4.98%  "
$ws.Range("E10").Value = "  +4.70%  "
$ws.Range("D11").Value = "0.07482"
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").Value = "0.9995"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "6.409"
$ws.Range("E13").Value = "  +6.02%  "
$ws.Range("D14").Value = "20.39"
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("D15").Value = "7.068"
$ws.Range("E15").Value = "  +7.27%  "
$ws.Range("D16").Value = "1.741.08"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").Value = "0.00001081"
$ws.Range("E17").Value = "  +3.49%  "
$ws.Range("D18").Value = "0.06707"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "82.70"
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  +6.19%  "
$ws.Range("D22").Value = "6.224"
$ws.Range("E22").Value = "  +5.61%  "
$ws.Range("D23").Value = "12.81"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "26.789.13"
$ws.Range("E24").Value = "  +7.76%  "
$ws.Range("D25").Value = "2.461"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "1.478"
$ws.Range("E26").Value = "  +25.86%  "
$ws.Range("D27").Value = "2.407"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Value = "152.59"
$ws.Range("E28").Value = "  +3.15%  "
$ws.Range("D29").Value = "19.64"
$ws.Range("E29").Value = "  +5.37%  "
$ws.Range("D30").Value = "1.936.75"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("D31").Value = "132.25"
$ws.Range("E31").Value = "  +5.60%  "
$ws.Range("D32").Value = "4.123"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "6.032"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("D35").Value = "1.691"
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("D36").Value = "12.87"
$ws.Range("E36").Value = "  +5.55%  "
$ws.Range("D37").Value = "5.431"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").Value = "0.02351"
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").Value = "0.2178"
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("D40").Value = "0.06273"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").Value = "8.490"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D42").Value = "1.225"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "0.6270"
$ws.Range("E43").Value = "  +5.87%  "
$ws.Range("D44").Value = "14.31"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("D45").Value = "0.9987"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "3.921"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").Value = "0.6064"
$ws.Range("E47").Value = "  +6.11%  "
$ws.Range("D48").Value = "129.24"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("D49").Value = "2.074"
$ws.Range("E49").Value = "  +6.19%  "
$ws.Range("D50").Value = "0.07226"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").Value = "77.98"
$ws.Range("E51").Value = "  +4.99%  "

# Restore the default (Normal) style so no residual number-format style
# is left attached to the Price cells.
$priceRange.Style = "Normal"

Write-Output "Updated cryptos list"
